# Regenerate merged AHB files
# 1) Rename the "_old" / "_new" header-name suffixes to "_FV2310" / "_FV2404"
# 2) Turn the used range into an Excel Table ("Table1")
# 3) Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fields = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $fields.Length; $i++) {
    $ws.Range($leftCols[$i]  + "1").Value = $fields[$i] + "_FV2310"
    $ws.Range($rightCols[$i] + "1").Value = $fields[$i] + "_FV2404"
}

# Build the table over the full used range A1:U59
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U59"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the top row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
